$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1405.7878
$ws.Range("J17").Value = 1104
$ws.Range("L17").Value = 3312
$ws.Range("N17").Value = -3648
$ws.Range("H51").Value = 5684.625
$ws.Range("I51").Value = 7120
$ws.Range("J51").Value = 4249.25
$ws.Range("K51").Value = 7120
$ws.Range("L51").Value = 4249.25
$ws.Range("M51").Value = -6636
$ws.Range("N51").Value = -5217.25
$ws.Range("H74").Value = 4198.1665
$ws.Range("J74").Value = 4198.3335
$ws.Range("L74").Value = 4198.3335
$ws.Range("N74").Value = -6070.3335
$ws.Range("H77").Value = 4198.1665
$ws.Range("J77").Value = 4198.3335
$ws.Range("L77").Value = 20991.6675
$ws.Range("N77").Value = -30351.6675
$ws.Range("H98").Value = 3377.6667
$ws.Range("I98").Value = 2742.7144
$ws.Range("K98").Value = 2742.7144
$ws.Range("M98").Value = -1244.7144
$ws.Range("H122").Value = 3377.6667
$ws.Range("I122").Value = 2742.7144
$ws.Range("K122").Value = 8228.143199999999
$ws.Range("M122").Value = -5778.143199999999
$ws.Range("H132").Value = 1248.5428
$ws.Range("I132").Value = 1142.3939
$ws.Range("K132").Value = 3427.1817
$ws.Range("M132").Value = -897.1817000000001
$ws.Range("H136").Value = 50000
$ws.Range("J136").Value = 50000
$ws.Range("L136").Value = 50000
$ws.Range("N136").Value = -60200
$ws.Range("H137").Value = 1650.6666
$ws.Range("I137").Value = 1223.375
$ws.Range("K137").Value = 3670.125
$ws.Range("M137").Value = -1120.125

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2735.4026
$ws.Range("I32").Value = 1608.409
$ws.Range("K32").Value = 1608.409
$ws.Range("M32").Value = -1321.409
$ws.Range("H132").Value = 1268.5454
$ws.Range("I132").Value = 1091.0541
$ws.Range("J132").Value = 2206.7144
$ws.Range("K132").Value = 3273.1623
$ws.Range("L132").Value = 6620.1432
$ws.Range("M132").Value = -743.1623
$ws.Range("N132").Value = -11680.1432

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2166.2307
$ws.Range("I107").Value = 1659.375
$ws.Range("J107").Value = 2977.2
$ws.Range("K107").Value = 1659.375
$ws.Range("L107").Value = 2977.2
$ws.Range("M107").Value = 260.625
$ws.Range("N107").Value = -6817.2

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 694.75
$ws.Range("I22").Value = 259.66666
$ws.Range("K22").Value = 259.66666
$ws.Range("M22").Value = 90.33334000000002
$ws.Range("H31").Value = 2940.0356
$ws.Range("I31").Value = 1205.6666
$ws.Range("K31").Value = 1205.6666
$ws.Range("M31").Value = -910.6666
$ws.Range("H34").Value = 2940.0356
$ws.Range("I34").Value = 1205.6666
$ws.Range("K34").Value = 1205.6666
$ws.Range("M34").Value = -1003.6666
$ws.Range("H109").Value = 50000
$ws.Range("J109").Value = 50000
$ws.Range("L109").Value = 50000
$ws.Range("N109").Value = -52080
$ws.Range("H122").Value = 1566.9062
$ws.Range("I122").Value = 1406.95
$ws.Range("J122").Value = 1833.5
$ws.Range("K122").Value = 4220.85
$ws.Range("L122").Value = 5500.5
$ws.Range("M122").Value = -1770.85
$ws.Range("N122").Value = -10400.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 91.111115
$ws.Range("I12").Value = 52.333332
$ws.Range("K12").Value = 156.999996
$ws.Range("M12").Value = 16.00000399999999
$ws.Range("H87").Value = 11390.2
$ws.Range("I87").Value = 2317
$ws.Range("K87").Value = 6951
$ws.Range("M87").Value = -5703
$ws.Range("H90").Value = 11390.2
$ws.Range("I90").Value = 2317
$ws.Range("K90").Value = 20853
$ws.Range("M90").Value = -14613
$ws.Range("H98").Value = 748.5
$ws.Range("I98").Value = 700
$ws.Range("J98").Value = 797
$ws.Range("K98").Value = 2100
$ws.Range("L98").Value = 2391
$ws.Range("M98").Value = -602
$ws.Range("N98").Value = -5387
$ws.Range("H103").Value = 1683
$ws.Range("I103").Value = 1099.8
$ws.Range("J103").Value = 4599
$ws.Range("K103").Value = 3299.4
$ws.Range("L103").Value = 13797
$ws.Range("M103").Value = -2420.4
$ws.Range("N103").Value = -15555
$ws.Range("H107").Value = 573.7059
$ws.Range("J107").Value = 573.7059
$ws.Range("L107").Value = 1721.1177
$ws.Range("N107").Value = -5561.117700000001
$ws.Range("H131").Value = 8635203
$ws.Range("I131").Value = 71429160
$ws.Range("J131").Value = 16424.275
$ws.Range("K131").Value = 214287480
$ws.Range("L131").Value = 49272.825
$ws.Range("M131").Value = -214282440
$ws.Range("N131").Value = -59352.825
$ws.Range("H132").Value = 1643.3572
$ws.Range("J132").Value = 1852.5
$ws.Range("L132").Value = 16672.5
$ws.Range("N132").Value = -21732.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1840.9131
$ws.Range("I102").Value = 1728.9166
$ws.Range("K102").Value = 1728.9166
$ws.Range("M102").Value = -106.9166
$ws.Range("H107").Value = 802.2
$ws.Range("I107").Value = 400.66666
$ws.Range("J107").Value = 974.2857
$ws.Range("K107").Value = 400.66666
$ws.Range("L107").Value = 974.2857
$ws.Range("M107").Value = 1519.33334
$ws.Range("N107").Value = -4814.2857
$ws.Range("H126").Value = 65762.06
$ws.Range("I126").Value = 3512.7856
$ws.Range("K126").Value = 10538.3568
$ws.Range("M126").Value = -8068.356800000001
$ws.Range("H132").Value = 3062.7837
$ws.Range("I132").Value = 2509.4333
$ws.Range("J132").Value = 5434.2856
$ws.Range("K132").Value = 7528.2999
$ws.Range("L132").Value = 16302.8568
$ws.Range("M132").Value = -4998.2999
$ws.Range("N132").Value = -21362.8568

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 10190.286
$ws.Range("I40").Value = 7203.3335
$ws.Range("J40").Value = 12430.5
$ws.Range("K40").Value = 7203.3335
$ws.Range("L40").Value = 12430.5
$ws.Range("M40").Value = -7067.3335
$ws.Range("N40").Value = -12702.5
$ws.Range("H50").Value = 27000
$ws.Range("J50").Value = 27000
$ws.Range("L50").Value = 27000
$ws.Range("N50").Value = -28274
$ws.Range("H61").Value = 3416
$ws.Range("I61").Value = 2963.4285
$ws.Range("K61").Value = 2963.4285
$ws.Range("M61").Value = -2761.4285
$ws.Range("H113").Value = 3416
$ws.Range("I113").Value = 2963.4285
$ws.Range("K113").Value = 2963.4285
$ws.Range("M113").Value = -793.4285
$ws.Range("H122").Value = 6537.36
$ws.Range("I122").Value = 4786.846
$ws.Range("K122").Value = 14360.538
$ws.Range("M122").Value = -11910.538

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("N26").ClearContents()
$ws.Range("H29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()
$ws.Range("H41").Value = 13374.75
$ws.Range("J41").Value = 13374.75
$ws.Range("L41").Value = 13374.75
$ws.Range("N41").Value = -14154.75
$ws.Range("H122").Value = 20333.963
$ws.Range("I122").Value = 32977.5
$ws.Range("J122").Value = 1943.3636
$ws.Range("K122").Value = 98932.5
$ws.Range("L122").Value = 5830.0908
$ws.Range("M122").Value = -96482.5
$ws.Range("N122").Value = -10730.0908
